$d = $word.ActiveDocument

function Find-InRange($rangeObj, [string]$text) {
    return $rangeObj.Find.Execute($text, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
}

function Replace-TextInParagraph($paraIndex, [string]$searchText, [string]$newText) {
    $p = $d.Paragraphs($paraIndex)
    $r = $d.Range($p.Range.Start, $p.Range.End)
    $found = Find-InRange $r $searchText
    if (-not $found) {
        throw "Could not find '$searchText' in paragraph $paraIndex"
    }
    $r.Text = $newText
}

# ===========================================================================
# Paragraph 5: "關鍵詞 (Biblica) (Chinese (Traditional)) is based on: Biblica
# Bible Dictionary, Biblica, Inc., 2023, which is licensed under a CC BY-SA
# 4.0 license."
#   -> "Biblica Study Notes (Key Terms) © 2023 Biblica Inc. Released under
#       CC BY-SA 4.0 license. Biblica Study Notes has been adapted in the
#       following languages: ... by Mission Mutual."
# ===========================================================================

# Step 1: remove "Biblica Bible Dictionary" through the hyperlinks and the
# trailing "." while the surrounding text is still unambiguous (this text
# appears only once in the whole document).
$p5 = $d.Paragraphs(5)
$rDictStart = $d.Range($p5.Range.Start, $p5.Range.End)
Find-InRange $rDictStart "Biblica Bible Dictionary" | Out-Null
$delStart = $rDictStart.Start

$p5 = $d.Paragraphs(5)
$rLink2 = $d.Range($p5.Range.Start, $p5.Range.End)
Find-InRange $rLink2 "CC BY-SA 4.0 license" | Out-Null
$rTrailingDot = $d.Range($rLink2.End, $p5.Range.End)
Find-InRange $rTrailingDot "." | Out-Null
$delEnd = $rTrailingDot.End

$d.Range($delStart, $delEnd).Delete() | Out-Null

# Step 2: update the remaining three runs (now unambiguous / safe to find
# fresh each time).
Replace-TextInParagraph 5 "關鍵詞 (Biblica)" "Biblica Study Notes (Key Terms)"
Replace-TextInParagraph 5 " (Chinese (Traditional)) is based on" " © 2023 Biblica Inc. Released under CC BY-SA 4.0 license. "
Replace-TextInParagraph 5 ": " "Biblica Study Notes"

# Step 3: append the long "has been adapted..." sentence right before the
# paragraph mark that ends paragraph 5.
$adaptedText = " has been adapted in the following languages: Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文)from Biblica Study Notes © 2023 Biblica Inc. Released under CC BY-SA 4.0 license by Mission Mutual."
$p5 = $d.Paragraphs(5)
$insertionPoint = $d.Range($p5.Range.End - 1, $p5.Range.End - 1)
$insertionPoint.InsertAfter($adaptedText)

# ===========================================================================
# Remove the "This PDF version is provided under the same license."
# paragraph (now paragraph 6) entirely.
# ===========================================================================
$d.Paragraphs(6).Range.Delete() | Out-Null

# ===========================================================================
# Remove the "License Information" heading paragraph (paragraph 4) entirely.
# ===========================================================================
$d.Paragraphs(4).Range.Delete() | Out-Null

# ===========================================================================
# Remove the italic "救恩, 救主, 舊約聖經" paragraph entirely.
# ===========================================================================
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*救恩, 救主, 舊約聖經*") {
        $p.Range.Delete() | Out-Null
        break
    }
}
